$d = $word.ActiveDocument

$replacements = @(
    @("50×40=", "81×47="),
    @("43×19=", "46×28="),
    @("73×15=", "47×17="),
    @("96×93=", "18×20="),
    @("91×90=", "56×81="),
    @("98×67=", "45×67="),
    @("36×36=", "73×34="),
    @("87×48=", "52×17="),
    @("60×79=", "29×55="),
    @("16×61=", "85×27="),
    @("68×13=", "89×82="),
    @("29×69=", "84×80="),
    @("50×94=", "28×61="),
    @("79×77=", "80×65="),
    @("59×75=", "90×55="),
    @("67×51=", "55×79="),
    @("94×91=", "14×17="),
    @("93×15=", "32×42="),
    @("70×79=", "11×86="),
    @("93×64=", "12×92="),
    @("30×27=", "87×43="),
    @("69×66=", "23×28="),
    @("43×71=", "76×87="),
    @("65×18=", "64×79="),
    @("64×75=", "69×73=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
